$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to stay text so numeric-looking values
# (e.g. "235.41") are not silently reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '37.331.75'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '2.066.90'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '235.41'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '57.37'
$ws.Range("E8").Value = '  -1.23%  '
$ws.Range("E9").Value = '  +3.38%  '
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("D12").Value = '2.370.56'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").Value = '14.38'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("D14").Value = '20.71'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").Value = '0.774'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '5.17'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").Value = '2.066.46'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").Value = '37.291.24'
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = '69.57'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").Value = '226.72'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +2.04%  '
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").Value = '166.90'
$ws.Range("D27").Value = '8.92'
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").Value = '1.40'
$ws.Range("E28").Value = '  -6.19%  '
$ws.Range("D29").Value = '0.128'
$ws.Range("E29").Value = '  -0.43%  '
$ws.Range("E30").Value = '  -0.33%  '
$ws.Range("E31").Value = '  -0.91%  '
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").Value = '4.54'
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").Value = '2.47'
$ws.Range("E35").Value = '  -3.31%  '
$ws.Range("D36").Value = '1.79'
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").Value = '3.33'
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("E39").Value = '  -4.63%  '
$ws.Range("E40").Value = '  -0.87%  '
$ws.Range("D41").Value = '0.0958'
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("D42").Value = '97.61'
$ws.Range("E42").Value = '  +0.73%  '
$ws.Range("D43").Value = '1.480.40'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("E44").Value = '  +0.89%  '
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("D46").Value = '4.08'
$ws.Range("E46").Value = '  -10.10%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = '15.25'
$ws.Range("E48").Value = '  -4.23%  '
$ws.Range("D49").Value = '7.20'
$ws.Range("E49").Value = '  -0.31%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").Value = '2.257.28'
$ws.Range("E51").Value = '  +0.30%  '

# Restore default cell style (drop the temporary text-format override)
# while keeping the values stored as text.
$ws.Range("D2:E51").Style = "Normal"

